$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 244, pushing the existing rows
# (old 244-256) down to 246-258.
$ws.Rows.Item(244).Resize(2).Insert()

# --- New row 244: week of 2022-07-11, quality "Especial" ---
$ws.Range("A244").Value = 10
$ws.Range("B244").Value = "Vega Modelo de Temuco"
$ws.Range("C244").Value = "La Araucanía"
$ws.Range("D244").Value = 44753
$ws.Range("E244").Value = 9
$ws.Range("F244").Value = "Fruta"
$ws.Range("G244").Value = 100102
$ws.Range("H244").Value = "Cítricos"
$ws.Range("I244").Value = 100102006
$ws.Range("J244").Value = "Pomelo"
$ws.Range("K244").Value = "Start Ruby"
$ws.Range("L244").Value = "Especial"
$ws.Range("M244").Value = 95
$ws.Range("N244").Value = 16000
$ws.Range("O244").Value = 16000
$ws.Range("P244").Value = 16000
$ws.Range("Q244").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R244").Value = "Región de O'Higgins"
$ws.Range("S244").Value = 1143
$ws.Range("T244").Value = 14

# --- New row 245: week of 2022-07-11, quality "Primera" ---
$ws.Range("A245").Value = 10
$ws.Range("B245").Value = "Vega Modelo de Temuco"
$ws.Range("C245").Value = "La Araucanía"
$ws.Range("D245").Value = 44753
$ws.Range("E245").Value = 9
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100102
$ws.Range("H245").Value = "Cítricos"
$ws.Range("I245").Value = 100102006
$ws.Range("J245").Value = "Pomelo"
$ws.Range("K245").Value = "Start Ruby"
$ws.Range("L245").Value = "Primera"
$ws.Range("M245").Value = 65
$ws.Range("N245").Value = 12000
$ws.Range("O245").Value = 12000
$ws.Range("P245").Value = 12000
$ws.Range("Q245").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R245").Value = "Región de O'Higgins"
$ws.Range("S245").Value = 800
$ws.Range("T245").Value = 15

# Make sure D244/D245 carry the same date-formatted style as the rest of
# column D (style index 2 in the original workbook).
$ws.Range("D244:D245").NumberFormat = $ws.Range("D246").NumberFormat

Write-Host "Done"
